$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: the "Value" property is parameterized in this COM model, so plain
# ".Value" returns a method descriptor instead of invoking the getter.
# Use ".Value2" for both reads and writes instead.

# --- Rows 2 and 3: swap A (Id), Q (Ost), R (Nord); round Q/R to integers ---
$a2 = $ws.Range("A2").Value2
$a3 = $ws.Range("A3").Value2
$ws.Range("A2").Value2 = $a3
$ws.Range("A3").Value2 = $a2

$q2 = [double]$ws.Range("Q2").Value2
$q3 = [double]$ws.Range("Q3").Value2
$ws.Range("Q2").Value2 = [Math]::Round($q3)
$ws.Range("Q3").Value2 = [Math]::Round($q2)

$r2 = [double]$ws.Range("R2").Value2
$r3 = [double]$ws.Range("R3").Value2
$ws.Range("R2").Value2 = [Math]::Round($r3)
$ws.Range("R3").Value2 = [Math]::Round($r2)

# --- Rows 4 and 5: swap whole species records ---
$cols = @("A","B","D","E","F","G","H")
foreach ($col in $cols) {
    $v4 = $ws.Range("$col`4").Value2
    $v5 = $ws.Range("$col`5").Value2
    $ws.Range("$col`4").Value2 = $v5
    $ws.Range("$col`5").Value2 = $v4
}

$q4 = [double]$ws.Range("Q4").Value2
$q5 = [double]$ws.Range("Q5").Value2
$ws.Range("Q4").Value2 = [Math]::Round($q5)
$ws.Range("Q5").Value2 = [Math]::Round($q4)

$r4 = [double]$ws.Range("R4").Value2
$r5 = [double]$ws.Range("R5").Value2
$ws.Range("R4").Value2 = [Math]::Round($r5)
$ws.Range("R5").Value2 = [Math]::Round($r4)

# "Publik kommentar" (AC) moves from row 5 to row 4.
$ws.Range("AC4").Value2 = "På murken låga"
$ws.Range("AC5").ClearContents()

# Clear the Starttid (Z) and Sluttid (AB) columns for all data rows;
# these "00:00" placeholder values are removed entirely.
$ws.Range("Z2:Z5").ClearContents()
$ws.Range("AB2:AB5").ClearContents()
